# CIERRE 30 ABR 22
# Close out payroll receipts: move the week label from week 15
# (11-17 Apr 2022) to week 17 (25 Apr - 1 May 2022) and update the
# corresponding payment amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 holds the literal week-label text; H9, B27, H27 and B43 all reference
# it (directly or transitively) via formulas, so they recalculate on their
# own once B9 changes.
$ws.Range("B9").Value = "SEMANA   17  DEL    25      Al   01   DE   MAYO          2022"

# Updated payment figures for this closing.
$ws.Range("K21").Value = 1400
$ws.Range("E23").Value = 0
$ws.Range("E40").Value = 1250

# Leave the selection where the author left it when saving.
$ws.Range("E41").Select()

$wb.Application.Calculate()
